$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "YEAR"
$ws.Range("D1").Value = "Room and Board"
$ws.Range("E1").Value = "Other Expenses"

$ws.Range("C10").Select()
